$d = $word.ActiveDocument

# --- 1) Split the "(Jer 29,11-14a)" run into "(", "Jer" (wrapped in
#        spellcheck proofErr markers) and " 29,11-14a)" ---------------
$findRng = $d.Content
$found = $findRng.Find.Execute("(Jer 29,11-14a)")
if (-not $found) {
    throw "Could not find the '(Jer 29,11-14a)' text to split"
}

# Re-seat the found hit into a fresh Range: InsertXML on the Find's own
# receiver range silently appends instead of replacing in this host.
$targetRng = $d.Range($findRng.Start, $findRng.End)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/part" pkg:contentType="application/xml"><pkg:xmlData><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:i/><w:iCs/></w:rPr><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:i/><w:iCs/></w:rPr><w:t>Jer</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve"> 29,11-14a)</w:t></w:r></w:p></pkg:xmlData></pkg:part></pkg:package>'

$targetRng.InsertXML($xml)

# --- 2) Remove the trailing empty paragraph before the sectPr ---------
$lastPara = $d.Paragraphs.Last
$start = $lastPara.Range.Start
$end = $lastPara.Range.End
if ($lastPara.Range.Text -eq "`r" -and $start -gt 0) {
    # include the preceding paragraph mark so the empty paragraph
    # collapses away instead of merely clearing its (already empty) text
    $delRng = $d.Range($start - 1, $end)
    $delRng.Delete()
}
